# Update to align with ORBIT configs
#
# Refreshes the "Cumulative Capacity" series on each scenario sheet with
# recalculated ORBIT figures. Three sheets also gain/lose a trailing year
# row where the new capacity curve now ends a year earlier/later:
#   - "Baseline-Mid (CC)"  : row 26 (2055) removed
#   - "Moderate-Mid (SC)"  : row 24 (2053) removed
#   - "Expanded-High"      : row 20 (2049) added
$wb = $excel.ActiveWorkbook

# --- Baseline-Low ---
$ws = $wb.Worksheets.Item("Baseline-Low")
$ws.Range("B2").Value = 1677.930163750701
$ws.Range("B3").Value = 2326.82557742865
$ws.Range("B4").Value = 3071.442162774158
$ws.Range("B5").Value = 3813.186489436427
$ws.Range("B6").Value = 4471.069388088458
$ws.Range("B7").Value = 5240.165460762002
$ws.Range("B8").Value = 5991.025415111622
$ws.Range("B9").Value = 6738.872307978584
$ws.Range("B10").Value = 7506.872105826967
$ws.Range("B11").Value = 8364.110800979011
$ws.Range("B12").Value = 9306.30117920732
$ws.Range("B13").Value = 10277.69373079723
$ws.Range("B14").Value = 11267.92801500317
$ws.Range("B15").Value = 12340.31378137568
$ws.Range("B16").Value = 13376.62374306103
$ws.Range("B17").Value = 14556.33076731603
$ws.Range("B18").Value = 15551.27324557401
$ws.Range("B19").Value = 16541.54169241404
$ws.Range("B20").Value = 17531.30099946307
$ws.Range("B21").Value = 18553.29271574779
$ws.Range("B22").Value = 19609.45417951032
$ws.Range("B23").Value = 20778.7998683004
$ws.Range("B24").Value = 21841.76789504404
$ws.Range("B25").Value = 22937.34591886247
$ws.Range("B26").Value = 24177.53899456122

# --- Baseline-Mid (SC) ---
$ws = $wb.Worksheets.Item("Baseline-Mid (SC)")
$ws.Range("B2").Value = 1106.975238711766
$ws.Range("B3").Value = 1789.77180773588
$ws.Range("B4").Value = 2679.304511637657
$ws.Range("B5").Value = 3830.801807547724
$ws.Range("B6").Value = 5575.28616049561
$ws.Range("B7").Value = 7454.838220871127
$ws.Range("B8").Value = 9227.117625974732
$ws.Range("B9").Value = 11012.67880652528
$ws.Range("B10").Value = 12835.17561469777
$ws.Range("B11").Value = 14707.36400630029
$ws.Range("B12").Value = 16656.12081439872
$ws.Range("B13").Value = 18514.28368380099
$ws.Range("B14").Value = 20364.68695661856
$ws.Range("B15").Value = 22309.35690273768
$ws.Range("B16").Value = 23382.02691704847
$ws.Range("B17").Value = 24582.78133046177

# --- Baseline-Mid (CC) ---
$ws = $wb.Worksheets.Item("Baseline-Mid (CC)")
$ws.Range("B2").Value = 1106.975238711766
$ws.Range("B3").Value = 1789.77180773588
$ws.Range("B4").Value = 2317.614649165451
$ws.Range("B5").Value = 3108.410306271269
$ws.Range("B6").Value = 4132.305091068081
$ws.Range("B7").Value = 5295.178197064989
$ws.Range("B8").Value = 6350.082754054949
$ws.Range("B9").Value = 7819.997608619689
$ws.Range("B10").Value = 9262.444196831215
$ws.Range("B11").Value = 10708.63658873511
$ws.Range("B12").Value = 12159.62040651018
$ws.Range("B13").Value = 13610.06957485048
$ws.Range("B14").Value = 15063.41532861543
$ws.Range("B15").Value = 16513.37322415364
$ws.Range("B16").Value = 17958.97134747294
$ws.Range("B17").Value = 19608.25131226893
$ws.Range("B18").Value = 20615.72670305235
$ws.Range("B19").Value = 21205.80315066649
$ws.Range("B20").Value = 21797.15901440689
$ws.Range("B21").Value = 22380.67165662052
$ws.Range("B22").Value = 22960.91938795787
$ws.Range("B23").Value = 23564.77257016452
$ws.Range("B24").Value = 24167.35050783615
$ws.Range("B25").Value = 24756.57496468689
$ws.Rows.Item(26).Delete()

# --- Moderate-Low ---
$ws = $wb.Worksheets.Item("Moderate-Low")
$ws.Range("B2").Value = 1992.536905323735
$ws.Range("B3").Value = 3111.416267637318
$ws.Range("B4").Value = 4336.207263445299
$ws.Range("B5").Value = 5567.884476013608
$ws.Range("B6").Value = 6708.001298602515
$ws.Range("B7").Value = 7952.911249180831
$ws.Range("B8").Value = 9255.80789060182
$ws.Range("B9").Value = 11287.71256951713
$ws.Range("B10").Value = 13311.98782436536
$ws.Range("B11").Value = 15449.02604170104
$ws.Range("B12").Value = 17651.19698277635
$ws.Range("B13").Value = 19753.28364129357
$ws.Range("B14").Value = 21267.92801500317
$ws.Range("B15").Value = 22340.31378137568
$ws.Range("B16").Value = 23376.62374306103
$ws.Range("B17").Value = 24556.33076731603
$ws.Range("B18").Value = 25551.27324557401
$ws.Range("B19").Value = 26541.54169241404
$ws.Range("B20").Value = 27531.30099946307
$ws.Range("B21").Value = 28553.29271574779
$ws.Range("B22").Value = 29609.45417951032
$ws.Range("B23").Value = 30778.7998683004
$ws.Range("B24").Value = 31841.76789504404
$ws.Range("B25").Value = 32937.34591886247
$ws.Range("B26").Value = 34177.53899456122

# --- Moderate-Mid (SC) ---
$ws = $wb.Worksheets.Item("Moderate-Mid (SC)")
$ws.Range("B2").Value = 1106.975238711766
$ws.Range("B3").Value = 2104.378549308913
$ws.Range("B4").Value = 3463.895201846325
$ws.Range("B5").Value = 5095.566908218865
$ws.Range("B6").Value = 6968.925701767374
$ws.Range("B7").Value = 8969.653240774353
$ws.Range("B8").Value = 10856.96648750686
$ws.Range("B9").Value = 12814.49945473534
$ws.Range("B10").Value = 15507.38722164896
$ws.Range("B11").Value = 18172.49265951406
$ws.Range("B12").Value = 20821.34877787585
$ws.Range("B13").Value = 23446.01547704969
$ws.Range("B14").Value = 25967.79259285823
$ws.Range("B15").Value = 27968.67806990608
$ws.Range("B16").Value = 29450.36582451781
$ws.Range("B17").Value = 31139.41455231239
$ws.Range("B18").Value = 32045.5110107284
$ws.Range("B19").Value = 32540.14680971203
$ws.Range("B20").Value = 33036.5296803653
$ws.Range("B21").Value = 33536.5296803653
$ws.Range("B22").Value = 34036.5296803653
$ws.Range("B23").Value = 34536.5296803653
$ws.Rows.Item(24).Delete()

# --- Expanded-High ---
$ws = $wb.Worksheets.Item("Expanded-High")
$ws.Range("B2").Value = 1106.975238711766
$ws.Range("B3").Value = 1789.77180773588
$ws.Range("B4").Value = 3495.535124433279
$ws.Range("B5").Value = 5461.032894852641
$ws.Range("B6").Value = 7665.663761242734
$ws.Range("B7").Value = 10457.45146654837
$ws.Range("B8").Value = 13236.7149592529
$ws.Range("B9").Value = 16239.80859929818
$ws.Range("B10").Value = 19760.73244453446
$ws.Range("B11").Value = 23717.76752940902
$ws.Range("B12").Value = 27863.67516887463
$ws.Range("B13").Value = 32303.26250546046
$ws.Range("B14").Value = 36814.00426987503
$ws.Range("B15").Value = 41329.17008787573
$ws.Range("B16").Value = 45583.5963688373
$ws.Range("B17").Value = 49660.04491817566
$ws.Range("B18").Value = 52536.26117655997
$ws.Range("B19").Value = 54058.62054410805
$ws.Range("A20").Value = 2049
$ws.Range("B20").Value = 54696.65732959851
